# Hortaliza, Agrícola del Norte S.A. de Arica - Ají
# Insert two new weekly price rows (74 and 75) into the daily-price log,
# shifting the existing rows 74:140 down to 76:142.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the existing data block (old row 74
# onward shifts down by 2, carrying its formatting with it).
$ws.Rows("74:75").Insert()

# --- New row 74 --------------------------------------------------------
$ws.Range("A74").Value = 1
$ws.Range("B74").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C74").Value = "Arica y Parinacota"
$ws.Range("D74").Value = 45096
$ws.Range("E74").Value = 15
$ws.Range("F74").Value = 100112021
$ws.Range("G74").Value = "Ají"
$ws.Range("H74").Value = "Inferno"
$ws.Range("I74").Value = "Segunda"
$ws.Range("J74").Value = 140
$ws.Range("K74").Value = 8000
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = 8500
$ws.Range("N74").Value = "$/caja 15 kilos"
$ws.Range("O74").Value = "Región de Arica y Parinacota"
$ws.Range("P74").Value = 567
$ws.Range("Q74").Value = 15
$ws.Range("R74").Value = "Hortaliza"

# --- New row 75 --------------------------------------------------------
$ws.Range("A75").Value = 1
$ws.Range("B75").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C75").Value = "Arica y Parinacota"
$ws.Range("D75").Value = 45096
$ws.Range("E75").Value = 15
$ws.Range("F75").Value = 100112021
$ws.Range("G75").Value = "Ají"
$ws.Range("H75").Value = "Inferno"
$ws.Range("I75").Value = "Tercera"
$ws.Range("J75").Value = 150
$ws.Range("K75").Value = 6000
$ws.Range("L75").Value = 7000
$ws.Range("M75").Value = 6500
$ws.Range("N75").Value = "$/caja 15 kilos"
$ws.Range("O75").Value = "Región de Arica y Parinacota"
$ws.Range("P75").Value = 433
$ws.Range("Q75").Value = 15
$ws.Range("R75").Value = "Hortaliza"
